# Update "想去人数" (F column) values on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same rows) to reflect newly scraped
# counts.

$wb = $excel.ActiveWorkbook

# Row -> new value mapping for column F (shared by both sheets).
$updates = @{
    2  = 296
    4  = 10250
    6  = 935
    7  = 1274
    8  = 6634
    11 = 192
    13 = 3153
    16 = 632
    18 = 521
    20 = 53
    21 = 1599
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
